$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 5
$ws.Range("E3").Value = "2021-05-14 10:49:41"
$ws.Range("E2").Value = "2021-05-24 17:39:41"

$ws.Range("D9").Select()
